$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.786.87"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "3.771.75"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.16"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.98"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").Value = "3.774.15"
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("E10").Value = "  +3.67%  "

$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.51"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "4.397.31"
$ws.Range("E15").Value = "  -0.89%  "

$ws.Range("D16").Value = "3.766.40"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").Value = "68.822.50"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("E20").Value = "  -0.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.00"
$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.74"
$ws.Range("E22").Value = "  +12.57%  "

$ws.Range("E23").Value = "  -1.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.75"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").Value = "  -2.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.43"
$ws.Range("E27").Value = "  +0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.54"
$ws.Range("E30").Value = "  +3.82%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.98"
$ws.Range("E32").Value = "  +0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.21"
$ws.Range("E33").Value = "  -2.22%  "

$ws.Range("D34").Value = "3.913.50"
$ws.Range("E34").Value = "  -0.93%  "

$ws.Range("D35").Value = "3.701.81"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.328"
$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "441.25"
$ws.Range("E42").Value = "  -4.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.99"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.53"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.80"
$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("D49").Value = "2.832.87"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("E51").Value = "  +0.65%  "
